$wb = $excel.ActiveWorkbook

# Add the new worksheet "FirstBankingLinks" as the last sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "FirstBankingLinks"

# Populate the new sheet's data
$newSheet.Range("A1").Value = "Chase First Banking: a debit card for teens and kids, managed by parents"
$newSheet.Range("A2").Value = "Chase First Banking: a debit card for teens and kids, managed by parents"
$newSheet.Range("A3").Value = "Chase First Banking: a debit card for teens and kids, managed by parents"
$newSheet.Range("A4").Value = "FAQs | Chase First Banking: child-friendly bank account opened by parents"
$newSheet.Range("A5").Value = "Chase First Banking vs. Chase High School Checking student accounts | Chase"

$newSheet.Range("A1:A5").Select()

# Make the new sheet the active one
$newSheet.Activate()

# Clear prior scroll state on the CarouselLinksTitles sheet (keep its existing selection)
$carousel = $wb.Worksheets.Item("CarouselLinksTitles")
$carousel.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

$tabsCompare = $wb.Worksheets.Item("TabsCompareAccounts")
$tabsCompare.Activate()
$tabsCompare.Range("F2").Select()

$newSheet.Activate()
$newSheet.Range("A1:A5").Select()
